# Updated cryptos list values (Price column D, Volume(1h) column E)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.545.46'
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").Value = '1.640.29'
$ws.Range("E3").Value = '  +2.35%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.16'
$ws.Range("E5").Value = '  +1.58%  '

$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3772'
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.41'
$ws.Range("E8").Value = '  +1.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3659'
$ws.Range("E9").Value = '  +1.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.273'
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08194'
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.01'
$ws.Range("E13").Value = '  +1.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.654'
$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001282'
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.420'
$ws.Range("E16").Value = '  +0.31%  '

$ws.Range("D17").Value = '1.640.90'
$ws.Range("E17").Value = '  +2.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.80'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06930'
$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.27'
$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.572'
$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").Value = '23.553.22'
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.85'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.065'
$ws.Range("E25").Value = '  +1.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.425'
$ws.Range("E26").Value = '  +1.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.32'
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.44'
$ws.Range("E28").Value = '  +0.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.361'
$ws.Range("E29").Value = '  +2.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.72'
$ws.Range("E30").Value = '  +1.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.385'
$ws.Range("E31").Value = '  -1.49%  '

$ws.Range("D32").Value = '1.824.30'
$ws.Range("E32").Value = '  +2.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.813'
$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9747'
$ws.Range("E34").Value = '  -0.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02829'
$ws.Range("E35").Value = '  +3.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.36'
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07387'
$ws.Range("E37").Value = '  -2.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2556'
$ws.Range("E38").Value = '  +2.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.194'
$ws.Range("E39").Value = '  +0.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08897'
$ws.Range("E40").Value = '  +1.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.384'
$ws.Range("E41").Value = '  +1.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7122'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.53'
$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.26'
$ws.Range("E44").Value = '  +4.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6555'
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.349'
$ws.Range("E46").Value = '  +1.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.045'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07984'
$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.50'
$ws.Range("E50").Value = '  -2.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.213'
$ws.Range("E51").Value = '  +0.38%  '
